$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 104; this shifts existing rows 104-129 down to 105-130
$ws.Rows.Item(104).Insert()

# Populate the new row 104 with the new data record
$ws.Cells.Item(104, 1).Value = 4
$ws.Cells.Item(104, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(104, 3).Value = "Los Lagos"
$ws.Cells.Item(104, 4).Value = 44782
$ws.Cells.Item(104, 5).Value = 10
$ws.Cells.Item(104, 6).Value = 100112022
$ws.Cells.Item(104, 7).Value = "Arveja Verde"
$ws.Cells.Item(104, 8).Value = "Perfection"
$ws.Cells.Item(104, 9).Value = "Primera"
$ws.Cells.Item(104, 10).Value = 40
$ws.Cells.Item(104, 11).Value = 44000
$ws.Cells.Item(104, 12).Value = 44000
$ws.Cells.Item(104, 13).Value = 44000
$ws.Cells.Item(104, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(104, 15).Value = "Provincia de Huasco"
$ws.Cells.Item(104, 16).Value = 1760
$ws.Cells.Item(104, 17).Value = 25
$ws.Cells.Item(104, 18).Value = "Hortaliza"
